# Generate Report for Handoff
#
# Re-stamps the localization-status report with a new handoff run:
#   - old GUID db9544c3-097d-4851-8c48-ad1e18d2d645 -> 0cc67afd-d66e-41ff-adb2-c9debb2fbae4
#   - old content hash 39008548413ef918a6482ba07b02402a5403885f
#         -> a49787555504f6fe67039e572c6e3a6b591d6f2a
#   - refreshed handoff timestamps
#
# Updates the cell text on all three sheets (Overview, zh-cn, de-de) and
# keeps each cell's hyperlink display text in sync with the new file names,
# without touching the hyperlink targets (their rIds / URLs are unchanged).

$wb = $excel.ActiveWorkbook

$oldGuid = "db9544c3-097d-4851-8c48-ad1e18d2d645"
$newGuid = "0cc67afd-d66e-41ff-adb2-c9debb2fbae4"
$oldHash = "39008548413ef918a6482ba07b02402a5403885f"
$newHash = "a49787555504f6fe67039e572c6e3a6b591d6f2a"

$newMdName    = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

# ---------------------------------------------------------------------
# Overview sheet: A2 file name + hyperlink, D2 latest-handoff timestamp
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$addrOverviewMd = "https://github.com/OpenLocalizationTest/oltest/blob/499a2e667e4dce1184c044d27b17b4269cc75764/e2e/$oldGuid.md"

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = "2016-59-13 04:59:34"

# Hyperlinks.Delete() only clears the collection as a whole in this host,
# so rebuild every link on the sheet (there is only the one here) with the
# refreshed display text, keeping the original target + relationship slot.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $addrOverviewMd, [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet: A2 file name, B2 extension (unchanged), D2 xlf + E2 time
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$addrZhMd  = "https://github.com/OpenLocalizationTest/oltest/blob/499a2e667e4dce1184c044d27b17b4269cc75764/e2e/$oldGuid.md"
$addrZhExt = "https://github.com/OpenLocalizationTest/oltest/blob/499a2e667e4dce1184c044d27b17b4269cc75764/e2e/$oldGuid.md"
$addrZhXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/79d288cac00e8afcdddedaa00628656bcf5ab701/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlfName
$wsZh.Range("E2").Value = "2016-03-13 04:59:31"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $addrZhMd,  [Type]::Missing, [Type]::Missing, $newMdName)    | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $addrZhExt, [Type]::Missing, [Type]::Missing, ".md")         | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $addrZhXlf, [Type]::Missing, [Type]::Missing, $newZhXlfName) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: A2 file name, B2 extension (unchanged), D2 xlf + E2 time
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$addrDeMd  = "https://github.com/OpenLocalizationTest/oltest/blob/499a2e667e4dce1184c044d27b17b4269cc75764/e2e/$oldGuid.md"
$addrDeExt = "https://github.com/OpenLocalizationTest/oltest/blob/499a2e667e4dce1184c044d27b17b4269cc75764/e2e/$oldGuid.md"
$addrDeXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/49f133af9c4cb8c70fe47386da11e2646424aaad/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlfName
$wsDe.Range("E2").Value = "2016-03-13 04:59:34"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $addrDeMd,  [Type]::Missing, [Type]::Missing, $newMdName)    | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $addrDeExt, [Type]::Missing, [Type]::Missing, ".md")         | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $addrDeXlf, [Type]::Missing, [Type]::Missing, $newDeXlfName) | Out-Null
